# Weekly update: a new price record (week of 2022-08-03) is added to the
# "Macroferia Regional de Talca - Ajo" sheet. The new record is inserted
# as a new row 292, pushing all existing rows from 292..333 down to 293..334.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 292, shifting rows 292-333
# down to 293-334 (matches the diff: old row292->new row293, ..., old
# row333->new row334).
$ws.Rows.Item(292).Insert()

# Populate the newly inserted row 292 with this week's record.
$ws.Cells.Item(292, 1).Value  = 5
$ws.Cells.Item(292, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(292, 3).Value  = "Maule"
$ws.Cells.Item(292, 4).Value  = 44776
$ws.Cells.Item(292, 5).Value  = 7
$ws.Cells.Item(292, 6).Value  = 100112003
$ws.Cells.Item(292, 7).Value  = "Ajo"
$ws.Cells.Item(292, 8).Value  = "Chino"
$ws.Cells.Item(292, 9).Value  = "Primera"
$ws.Cells.Item(292, 10).Value = 300
$ws.Cells.Item(292, 11).Value = 28000
$ws.Cells.Item(292, 12).Value = 28000
$ws.Cells.Item(292, 13).Value = 28000
$ws.Cells.Item(292, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(292, 15).Value = "China"
$ws.Cells.Item(292, 16).Value = 2800
$ws.Cells.Item(292, 17).Value = 10
$ws.Cells.Item(292, 18).Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D.
$ws.Cells.Item(292, 4).NumberFormat = $ws.Cells.Item(293, 4).NumberFormat
